$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 (M2:T2) with new TPM-derived values
$ws.Range("M2").Value = 0.015498
$ws.Range("N2").Value = 0.046494
$ws.Range("O2").Value = 0.2487028339734469
$ws.Range("P2").Value = 0.2487028339734469
$ws.Range("Q2").Value = 0.009887651676
$ws.Range("R2").Value = 0.088988865084
$ws.Range("S2").Value = 0.2487028339734469
$ws.Range("T2").Value = 0.2487028339734469

# Update row 3 (O3,P3,S3,T3) with new TPM-derived values
$ws.Range("O3").Value = 0.62532496014892
$ws.Range("P3").Value = 0.62532496014892
$ws.Range("S3").Value = 0.62532496014892
$ws.Range("T3").Value = 0.62532496014892

# Add new row 4 (FAPs -> Tac1/Tacr1 -> MuSCs interaction)
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Tac1"
$ws.Range("C4").Value = "Tacr1"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.6379953333333334
$ws.Range("H4").Value = 1.913986
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.007850000000000001
$ws.Range("N4").Value = 0.02355
$ws.Range("O4").Value = 0.1259722058776331
$ws.Range("P4").Value = 0.1259722058776331
$ws.Range("Q4").Value = 0.005008263366666667
$ws.Range("R4").Value = 0.0450743703
$ws.Range("S4").Value = 0.1259722058776331
$ws.Range("T4").Value = 0.1259722058776331
